# Generate Report for Handback
# Update the timestamp strings recorded for the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file, and
# de-de's "Correspond Handoff Datetime" for the same file (shared text).
$wsOverview.Range("G2").Value = "2016-08-30 05:05:34"
$wsDeDe.Range("H2").Value = "2016-08-30 05:05:34"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-30 05:05:29"
$wsZhCn.Range("K2").Value = "2016-08-30 05:05:46"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-30 05:05:53"
